$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: title paragraph.
# The title text was originally typed as two separate runs with a
# "_GoBack" bookmark (Word's "last edit" marker) sitting between them:
#   "Group 1 - Crime Incident Repor" | _GoBack | "ts in Boston Between June 2015 - September 2018"
# It becomes a single contiguous run (bookmark gone from here - it will
# be relocated below) with the full title text.
# ---------------------------------------------------------------------
$titleOld = "Group 1 " + [char]0x2013 + " Crime Incident Report"
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Execute($titleOld, $true, $false, $false, $false, $false, $true, 1, $false, $titleOld, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: "... documents the initial incident between 2018-2018."
# becomes "... between 2015-2018." - the stray duplicate "2018" start
# year is corrected to "2015". The edit also leaves Word's "_GoBack"
# bookmark positioned right after the corrected digit (i.e. between
# "5" and "-2018"), matching where the author's cursor would land
# after typing the replacement digit.
# ---------------------------------------------------------------------
$rngFind = $d.Content.Duplicate
$find2 = $rngFind.Find
$find2.ClearFormatting()
$find2.Execute("2018-2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$spanStart = $rngFind.Start              # start of the first "2018"
$spanEnd = $rngFind.End                  # just past the second "2018" (before the ".")

$posAfter201 = $spanStart + 3            # between "201" and the stray "8"
$posAfterDigit = $spanStart + 4          # between the stray "8" and "-2018"

# Pin the run-split points with throwaway bookmarks first: plain text
# edits cause neighboring same-formatted runs in the paragraph to be
# rejoined, so pinning here keeps "  between 201" | "<digit>" | "-2018"
# | "." as four distinct runs once the character is changed below.
$d.Bookmarks.Add("zzEditPin1", $d.Range($posAfter201, $posAfter201)) | Out-Null
$d.Bookmarks.Add("zzEditPin2", $d.Range($posAfterDigit, $posAfterDigit)) | Out-Null
$d.Bookmarks.Add("zzEditPin3", $d.Range($spanEnd, $spanEnd)) | Out-Null

# Correct the stray "8" (making "2018-2018" read "2015-2018").
$digitRng = $d.Range($posAfter201, $posAfterDigit)
$digitRng.Text = "5"

# Move (recreate) the "_GoBack" bookmark to sit right after the digit
# we just fixed, then drop the scaffolding bookmarks.
$pin2Range = $d.Bookmarks("zzEditPin2").Range
$d.Bookmarks.Add("_GoBack", $pin2Range) | Out-Null

$d.Bookmarks("zzEditPin1").Delete()
$d.Bookmarks("zzEditPin2").Delete()
$d.Bookmarks("zzEditPin3").Delete()

Write-Output "Edits applied."
